$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E15").Value = "[]"

$ws.Range("D24").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E24").Value = "['Normal']"

$ws.Range("D25").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E25").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D38").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E38").Value = "['SoftwareFault']"

$ws.Range("D39").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E39").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

$ws.Range("D54").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E54").Value = "[]"

$ws.Range("D58").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E58").Value = "['Normal']"

$ws.Range("D61").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E61").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D67").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E67").Value = "['Normal', 'HardwareFault']"

$ws.Range("D68").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E68").Value = "['Normal', 'SurroundingEnvironment']"

$ws.Range("D73").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'ParamViolation']"

$ws.Range("D75").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E75").Value = "['SoftwareFault']"

$ws.Range("D84").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E84").Value = "['Normal', 'SurroundingEnvironment']"

$ws.Range("D109").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E109").Value = "['Normal', 'SurroundingEnvironment']"
